# Auto-generated market-data refresh for Spriggan_Profits workbook.
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) per-sheet
# with freshly pulled market values; mirrors the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 5163.8
$ws.Range("I9").Value = 6843.467
$ws.Range("J9").Value = 124.8
$ws.Range("K9").Value = 6843.467
$ws.Range("L9").Value = 124.8
$ws.Range("M9").Value = -6674.467
$ws.Range("N9").Value = -462.8
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 2624.1853
$ws.Range("J132").Value = 1249.5
$ws.Range("L132").Value = 3748.5
$ws.Range("N132").Value = -8808.5
$ws.Range("H137").Value = 4874.5586
$ws.Range("I137").Value = 4754.2
$ws.Range("K137").Value = 14262.6
$ws.Range("M137").Value = -11712.6
$ws.Range("H138").Value = 4817.4
$ws.Range("J138").Value = 5022.128
$ws.Range("L138").Value = 15066.384
$ws.Range("N138").Value = -25346.384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25628
$ws.Range("I32").Value = 19748.834
$ws.Range("J32").Value = 29547.445
$ws.Range("K32").Value = 19748.834
$ws.Range("L32").Value = 29547.445
$ws.Range("M32").Value = -19461.834
$ws.Range("N32").Value = -30121.445
$ws.Range("H44").Value = 49999
$ws.Range("J44").Value = 49999
$ws.Range("L44").Value = 49999
$ws.Range("N44").Value = -50975
$ws.Range("H45").Value = 1997.1666
$ws.Range("I45").Value = 1996.5
$ws.Range("K45").Value = 1996.5
$ws.Range("M45").Value = -1619.5
$ws.Range("H61").Value = 40820348
$ws.Range("J61").Value = 4858.0835
$ws.Range("L61").Value = 4858.0835
$ws.Range("N61").Value = -5282.0835
$ws.Range("H110").Value = 73475.86
$ws.Range("I110").Value = 112592.336
$ws.Range("J110").Value = 3066.2
$ws.Range("K110").Value = 112592.336
$ws.Range("L110").Value = 3066.2
$ws.Range("M110").Value = -110547.336
$ws.Range("N110").Value = -7156.2
$ws.Range("H122").Value = 7179.8
$ws.Range("I122").Value = 7179.8
$ws.Range("K122").Value = 21539.4
$ws.Range("M122").Value = -19089.4
$ws.Range("H132").Value = 2574386
$ws.Range("I132").Value = 2712190.5
$ws.Range("K132").Value = 8136571.5
$ws.Range("M132").Value = -8134041.5
$ws.Range("H136").Value = 40820348
$ws.Range("J136").Value = 4858.0835
$ws.Range("L136").Value = 14574.2505
$ws.Range("N136").Value = -19674.2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 76990
$ws.Range("J52").Value = 76990
$ws.Range("L52").Value = 76990
$ws.Range("N52").Value = -77516
$ws.Range("H121").Value = 76990
$ws.Range("J121").Value = 76990
$ws.Range("L121").Value = 76990
$ws.Range("N121").Value = -80484
$ws.Range("H134").Value = 11366434
$ws.Range("I134").Value = 11630583
$ws.Range("K134").Value = 34891749
$ws.Range("M134").Value = -34889214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2021504.6
$ws.Range("I16").Value = 5050755
$ws.Range("K16").Value = 5050755
$ws.Range("M16").Value = -5050468
$ws.Range("H31").Value = 8750.950999999999
$ws.Range("I31").Value = 5659.1816
$ws.Range("J31").Value = 12330.895
$ws.Range("K31").Value = 5659.1816
$ws.Range("L31").Value = 12330.895
$ws.Range("M31").Value = -5364.1816
$ws.Range("N31").Value = -12920.895
$ws.Range("H34").Value = 8750.950999999999
$ws.Range("I34").Value = 5659.1816
$ws.Range("J34").Value = 12330.895
$ws.Range("K34").Value = 5659.1816
$ws.Range("L34").Value = 12330.895
$ws.Range("M34").Value = -5457.1816
$ws.Range("N34").Value = -12734.895
$ws.Range("H99").Value = 3351.3
$ws.Range("I99").Value = 3617.6924
$ws.Range("J99").Value = 2856.5715
$ws.Range("K99").Value = 3617.6924
$ws.Range("L99").Value = 2856.5715
$ws.Range("M99").Value = -2119.6924
$ws.Range("N99").Value = -5852.5715
$ws.Range("H106").Value = 15750
$ws.Range("J106").Value = 15750
$ws.Range("L106").Value = 15750
$ws.Range("N106").Value = -18274
$ws.Range("H113").Value = 2021504.6
$ws.Range("I113").Value = 5050755
$ws.Range("K113").Value = 5050755
$ws.Range("M113").Value = -5048585
$ws.Range("H126").Value = 3351.3
$ws.Range("I126").Value = 3617.6924
$ws.Range("J126").Value = 2856.5715
$ws.Range("K126").Value = 10853.0772
$ws.Range("L126").Value = 8569.7145
$ws.Range("M126").Value = -8383.0772
$ws.Range("N126").Value = -13509.7145
$ws.Range("H134").Value = 6758052.5
$ws.Range("I134").Value = 7354252.5
$ws.Range("J134").Value = 1116.6666
$ws.Range("K134").Value = 22062757.5
$ws.Range("L134").Value = 3349.9998
$ws.Range("M134").Value = -22060222.5
$ws.Range("N134").Value = -8419.9998
$ws.Range("H141").Value = 276255.1
$ws.Range("J141").Value = 398107.94
$ws.Range("L141").Value = 398107.94
$ws.Range("N141").Value = -408467.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.736843
$ws.Range("J2").Value = 55.77778
$ws.Range("L2").Value = 334.66668
$ws.Range("N2").Value = -560.66668
$ws.Range("H4").Value = 763570.6
$ws.Range("I4").Value = 942116.6
$ws.Range("K4").Value = 2826349.8
$ws.Range("M4").Value = -2826237.8
$ws.Range("H9").Value = 4443
$ws.Range("J9").Value = 4443
$ws.Range("L9").Value = 13329
$ws.Range("N9").Value = -13777
$ws.Range("H32").Value = 60001210
$ws.Range("J32").Value = 322.5
$ws.Range("L32").Value = 967.5
$ws.Range("N32").Value = -1533.5
$ws.Range("H35").Value = 350
$ws.Range("I35").Value = 350
$ws.Range("K35").Value = 1050
$ws.Range("M35").Value = -762
$ws.Range("H39").Value = 2552.7058
$ws.Range("I39").Value = 1028.5
$ws.Range("J39").Value = 9665.666999999999
$ws.Range("K39").Value = 3085.5
$ws.Range("L39").Value = 28997.001
$ws.Range("M39").Value = -2791.5
$ws.Range("N39").Value = -29585.001
$ws.Range("H107").Value = 997.069
$ws.Range("I107").Value = 312.08334
$ws.Range("J107").Value = 1480.5883
$ws.Range("K107").Value = 936.2500200000001
$ws.Range("L107").Value = 4441.7649
$ws.Range("M107").Value = 983.7499799999999
$ws.Range("N107").Value = -8281.7649
$ws.Range("H113").Value = 201555.8
$ws.Range("I113").Value = 251452.5
$ws.Range("J113").Value = 1969
$ws.Range("K113").Value = 754357.5
$ws.Range("L113").Value = 5907
$ws.Range("M113").Value = -752187.5
$ws.Range("N113").Value = -10247
$ws.Range("H127").Value = 2999
$ws.Range("J127").Value = 2999
$ws.Range("L127").Value = 8997
$ws.Range("N127").Value = -18917
$ws.Range("H132").Value = 2560.4
$ws.Range("I132").Value = 2599
$ws.Range("J132").Value = 2502.5
$ws.Range("K132").Value = 23391
$ws.Range("L132").Value = 22522.5
$ws.Range("M132").Value = -20861
$ws.Range("N132").Value = -27582.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 847.375
$ws.Range("I97").Value = 528.1111
$ws.Range("J97").Value = 1257.8572
$ws.Range("K97").Value = 528.1111
$ws.Range("L97").Value = 1257.8572
$ws.Range("M97").Value = -32.11109999999996
$ws.Range("N97").Value = -2249.8572
$ws.Range("H126").Value = 10568.154
$ws.Range("I126").Value = 9708.143
$ws.Range("J126").Value = 11571.5
$ws.Range("K126").Value = 29124.429
$ws.Range("L126").Value = 34714.5
$ws.Range("M126").Value = -26654.429
$ws.Range("N126").Value = -39654.5
$ws.Range("H132").Value = 2911596.8
$ws.Range("I132").Value = 2980712
$ws.Range("K132").Value = 8942136
$ws.Range("M132").Value = -8939606

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 61342.25
$ws.Range("J104").Value = 61342.25
$ws.Range("L104").Value = 61342.25
$ws.Range("N104").Value = -68330.25
$ws.Range("H122").Value = 6894.1577
$ws.Range("I122").Value = 6999
$ws.Range("J122").Value = 6750
$ws.Range("K122").Value = 20997
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -18547
$ws.Range("N122").Value = -25150
$ws.Range("H132").Value = 29591434
$ws.Range("I132").Value = 32281138
$ws.Range("K132").Value = 96843414
$ws.Range("M132").Value = -96840884
$ws.Range("H136").Value = 2697.225
$ws.Range("I136").Value = 2570.375
$ws.Range("J136").Value = 2887.5
$ws.Range("K136").Value = 7711.125
$ws.Range("L136").Value = 8662.5
$ws.Range("M136").Value = -5161.125
$ws.Range("N136").Value = -13762.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 12499.667
$ws.Range("J19").Value = 12499.667
$ws.Range("L19").Value = 12499.667
$ws.Range("N19").Value = -12847.667
$ws.Range("H104").Value = 19950
$ws.Range("J104").Value = 19950
$ws.Range("L104").Value = 19950
$ws.Range("N104").Value = -26938
$ws.Range("H107").Value = 1089.7778
$ws.Range("I107").Value = 483.625
$ws.Range("K107").Value = 1450.875
$ws.Range("M107").Value = 469.125
$ws.Range("H115").Value = 39999
$ws.Range("J115").Value = 39999
$ws.Range("L115").Value = 39999
$ws.Range("N115").Value = -43133
$ws.Range("H122").Value = 7666.6665
$ws.Range("I122").Value = 7666.6665
$ws.Range("K122").Value = 22999.9995
$ws.Range("M122").Value = -20549.9995
$ws.Range("H126").Value = 4499.421
$ws.Range("I126").Value = 1311.25
$ws.Range("K126").Value = 3933.75
$ws.Range("M126").Value = -1463.75
$ws.Range("H132").Value = 12200129
$ws.Range("I132").Value = 15627732
$ws.Range("K132").Value = 46883196
$ws.Range("M132").Value = -46880666
$ws.Range("H136").Value = 17243182
$ws.Range("I136").Value = 21741066
$ws.Range("J136").Value = 1294.1666
$ws.Range("K136").Value = 65223198
$ws.Range("L136").Value = 3882.4998
$ws.Range("M136").Value = -65220648
$ws.Range("N136").Value = -8982.4998
